# Update scripts with new TPM
#
# The underlying NATMI TPM recomputation dropped the "Inflammatory-Mac"
# sending cluster entirely and recalculated the ligand/receptor
# expression statistics for the remaining clusters (ECs, FAPs) against
# the same three target clusters (ECs, FAPs, MuSCs).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the three "Inflammatory-Mac" rows (previously rows 8-10). This
# also shrinks the used range from A1:T10 down to A1:T7 and, because the
# "Inflammatory-Mac" shared string then has no remaining references, it
# is dropped from the shared-strings table on save.
$ws.Range("A8:T10").Delete() | Out-Null

# --- Row 2 : ECs -> Hcrt -> Hcrtr1 -> ECs -------------------------------
$ws.Range("E2").Value2 = 2
$ws.Range("F2").Value2 = 0.6666666666666666
$ws.Range("G2").Value2 = 0.6051576666666666
$ws.Range("H2").Value2 = 1.815473
$ws.Range("I2").Value2 = 0.575135406723878
$ws.Range("J2").Value2 = 0.575135406723878
$ws.Range("K2").Value2 = 3
$ws.Range("L2").Value2 = 1
$ws.Range("M2").Value2 = 0.6609660000000001
$ws.Range("N2").Value2 = 1.982898
$ws.Range("O2").Value2 = 0.8866694927077409
$ws.Range("P2").Value2 = 0.8866694927077408
$ws.Range("Q2").Value2 = 0.399988642306
$ws.Range("R2").Value2 = 3.599897780754
$ws.Range("S2").Value2 = 0.5099550193181211
$ws.Range("T2").Value2 = 0.509955019318121

# --- Row 3 : ECs -> Hcrt -> Hcrtr1 -> FAPs -------------------------------
$ws.Range("E3").Value2 = 2
$ws.Range("F3").Value2 = 0.6666666666666666
$ws.Range("G3").Value2 = 0.6051576666666666
$ws.Range("H3").Value2 = 1.815473
$ws.Range("I3").Value2 = 0.575135406723878
$ws.Range("J3").Value2 = 0.575135406723878
$ws.Range("O3").Value2 = 0.08546493741571064
$ws.Range("P3").Value2 = 0.08546493741571062
$ws.Range("Q3").Value2 = 0.0385543932241111
$ws.Range("R3").Value2 = 0.346989539017
$ws.Range("S3").Value2 = 0.04915391154121551
$ws.Range("T3").Value2 = 0.0491539115412155

# --- Row 4 : ECs -> Hcrt -> Hcrtr1 -> MuSCs ------------------------------
$ws.Range("E4").Value2 = 2
$ws.Range("F4").Value2 = 0.6666666666666666
$ws.Range("G4").Value2 = 0.6051576666666666
$ws.Range("H4").Value2 = 1.815473
$ws.Range("I4").Value2 = 0.575135406723878
$ws.Range("J4").Value2 = 0.575135406723878
$ws.Range("K4").Value2 = 1
$ws.Range("L4").Value2 = 0.3333333333333333
$ws.Range("M4").Value2 = 0.02077233333333333
$ws.Range("N4").Value2 = 0.062317
$ws.Range("O4").Value2 = 0.02786556987654851
$ws.Range("P4").Value2 = 0.0278655698765485
$ws.Range("Q4").Value2 = 0.01257053677122222
$ws.Range("R4").Value2 = 0.113134830941
$ws.Range("S4").Value2 = 0.01602647586454137
$ws.Range("T4").Value2 = 0.01602647586454137

# --- Row 5 : FAPs -> Hcrt -> Hcrtr1 -> ECs -------------------------------
$ws.Range("I5").Value2 = 0.424864593276122
$ws.Range("J5").Value2 = 0.424864593276122
$ws.Range("K5").Value2 = 3
$ws.Range("L5").Value2 = 1
$ws.Range("M5").Value2 = 0.6609660000000001
$ws.Range("N5").Value2 = 1.982898
$ws.Range("O5").Value2 = 0.8866694927077409
$ws.Range("P5").Value2 = 0.8866694927077408
$ws.Range("Q5").Value2 = 0.295480003216
$ws.Range("R5").Value2 = 2.659320028944
$ws.Range("S5").Value2 = 0.3767144733896198
$ws.Range("T5").Value2 = 0.3767144733896197

# --- Row 6 : FAPs -> Hcrt -> Hcrtr1 -> FAPs ------------------------------
$ws.Range("I6").Value2 = 0.424864593276122
$ws.Range("J6").Value2 = 0.424864593276122
$ws.Range("O6").Value2 = 0.08546493741571064
$ws.Range("P6").Value2 = 0.08546493741571062
$ws.Range("S6").Value2 = 0.03631102587449513
$ws.Range("T6").Value2 = 0.03631102587449511

# --- Row 7 : FAPs -> Hcrt -> Hcrtr1 -> MuSCs -----------------------------
$ws.Range("I7").Value2 = 0.424864593276122
$ws.Range("J7").Value2 = 0.424864593276122
$ws.Range("K7").Value2 = 1
$ws.Range("L7").Value2 = 0.3333333333333333
$ws.Range("M7").Value2 = 0.02077233333333333
$ws.Range("N7").Value2 = 0.062317
$ws.Range("O7").Value2 = 0.02786556987654851
$ws.Range("P7").Value2 = 0.0278655698765485
$ws.Range("Q7").Value2 = 0.009286119286222222
$ws.Range("R7").Value2 = 0.08357507357599998
$ws.Range("S7").Value2 = 0.01183909401200714
$ws.Range("T7").Value2 = 0.01183909401200714
